$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for row 3: B3 gets a task description, C3 gets the "Will Do" note with wrapped text.
$ws.Range("B3").Value = "Design DB diagram"
$ws.Range("C3").Value = "- Analyze product backlog`n- Study about report"
$ws.Range("C3").WrapText = $true

# Move the active selection to D3, matching the saved view state.
$ws.Range("D3").Select() | Out-Null
